# "Huge v2 update" translation patch: the English translation column
# moves from D to B (so translations sit right next to the Japanese
# source text in column A), and the previously-untranslated row 2
# ("startTopFix" - an internal marker, not real dialogue) now gets a
# translation equal to its own key.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 14

# Snapshot existing column D translations before we start overwriting
# anything, then rebuild column B from them.
$translations = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $translations[$r] = $ws.Cells.Item($r, 4).Value2
}

for ($r = 1; $r -le $lastRow; $r++) {
    if ($r -eq 2) {
        # Row 2 had no translation before; give it the same value as
        # column A (the marker key itself).
        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
    } elseif (-not [string]::IsNullOrEmpty($translations[$r])) {
        $ws.Cells.Item($r, 2).Value = $translations[$r]
    }
}

# Column D is no longer used now that translations live in column B.
$ws.Range("D1:D$lastRow").ClearContents()
